$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ANG YI LING"
$ws.Range("B2").Value = "014-6263882"
$ws.Range("C2").Value = "yi_ling13@hotmail.com"
$ws.Range("D2").Value = "No"

$ws.Range("F2").Value = "[{'job_title': 'SENIOR CORPORATE TRAVEL CONSULTANT', 'job_company': 'FCM TRAVEL SOLUTIONS MALAYSIA', 'Industries': ['Travel and Tourism'], 'start_date': '2018-11', 'end_date': '2021-07', 'job_location': 'Malaysia', 'job_duration': '2 years 9 months'}, {'job_title': 'SENIOR CORPORATE TRAVEL CONSULTANT', 'job_company': 'FCM TRAVEL SINGAPORE', 'Industries': ['Travel and Tourism'], 'start_date': '2021-08', 'end_date': '2021-11', 'job_location': 'Singapore', 'job_duration': '3 months'}, {'job_title': 'CORPORATE TRAVEL CONSULTANT', 'job_company': 'HOLIDAY TOURS SDN BHD', 'Industries': ['Travel and Tourism'], 'start_date': '2016-07', 'end_date': '2018-04', 'job_location': 'Kuala Lumpur', 'job_duration': '1 year 9 months'}]"

$ws.Range("H2").Value = "[{'field_of_study': 'Tourism Management', 'level': ""Bachelor's Degree"", 'cgpa': '3.38', 'university': 'TUNKU ABDUL RAHMAN UNIVERSITY COLLEGE', 'start_date': '2016', 'year_of_graduation': '2016'}, {'field_of_study': 'Hospitality Management', 'level': 'Diploma', 'cgpa': '3.20', 'university': 'N/A', 'start_date': '2014', 'year_of_graduation': '2014'}]"

$ws.Range("I2").Value = "['Googe Analytics for Beginners', 'The Fundamental SQL Bootcamp', 'Python Programming for Beginners', 'Mastering SQL server']"

$ws.Range("J2").Value = "['Written and verbal communications in Chinese, English & Malay', 'Familiar with airlines reservation system, Sabre & Amadeus', 'Experienced with hotel distribution channels', 'Quick learner, ability to learn new skills quickly, act on feedback constructively and apply new knowledge immediately with the ability to identify learning opportunities']"

$ws.Range("K2").Value = "['Chinese', 'English', 'Malay']"

$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 0
$ws.Range("P2").Value = 0

$ws.Range("W2").Value = "2016, 2014, "

$ws.Range("Z2").Value = "**Alignments:**
1. **Education Background:**
   - The candidate has a Bachelor's Degree in Tourism Management and a Diploma in Hospitality Management, which aligns with the job description's requirement for a background in Tourism or Hospitality.
2. **Skills:**
   - The candidate has experience with airlines reservation systems like Sabre & Amadeus, which aligns with the job requirement for familiarity with reservation systems.
   - The candidate is proficient in written and verbal communications in Chinese, English, and Malay, which aligns with the job requirement for multilingual communication skills.
3. **Previous Job Roles:**
   - The candidate has held various roles in the Travel and Tourism industry, including Senior Corporate Travel Consultant positions, which align with the job description's requirement for relevant industry experience.
4. **Professional Certificates:**
   - The candidate has professional certificates in Google Analytics, SQL, and Python programming, which may align with the job requirement for technical skills or certifications.
**Misalignments:**
1. **Education Background:**
   - The candidate's CGPA for the Bachelor's Degree and Diploma is not mentioned in the provided information, which could be a misalignment if the job description requires a specific GPA.
2. **Skills:**
   - While the candidate mentions being a quick learner and adaptable, specific skills or experiences related to the job description may not be explicitly mentioned.
3. **Previous Job Roles:**
   - The candidate's previous job roles focus on travel consultancy, but specific experience or achievements related to the job description's requirements may not be highlighted.
4. **Professional Certificates:**
   - The relevance of the professional certificates to the job description's requirements is not explicitly stated, so there may be a misalignment in terms of specific certification needs.
Overall, the candidate's background in Tourism and Hospitality, relevant industry experience, and language skills align well with the job description. However, more explicit details on certain skills and experiences related to the job requirements could further strengthen the alignment."

$ws.Range("AA2").Value = 13
